$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.561916
$ws.Range("H2").Value = 4.685748
$ws.Range("I2").Value = 0.02125955722027638
$ws.Range("J2").Value = 0.02125955722027638
$ws.Range("M2").Value = 2.733663333333333
$ws.Range("N2").Value = 8.200989999999999
$ws.Range("O2").Value = 0.04037266183309663
$ws.Range("P2").Value = 0.04037266183309663
$ws.Range("Q2").Value = 4.269752498946667
$ws.Range("R2").Value = 38.42777249052
$ws.Range("S2").Value = 0.0008583049143755859
$ws.Range("T2").Value = 0.0008583049143755859
$ws.Range("G3").Value = 1.561916
$ws.Range("H3").Value = 4.685748
$ws.Range("I3").Value = 0.02125955722027638
$ws.Range("J3").Value = 0.02125955722027638
$ws.Range("O3").Value = 0.6389522306252696
$ws.Range("P3").Value = 0.6389522306252696
$ws.Range("Q3").Value = 67.57463490265333
$ws.Range("R3").Value = 608.1717141238801
$ws.Range("S3").Value = 0.01358384150800115
$ws.Range("T3").Value = 0.01358384150800115
$ws.Range("G4").Value = 1.561916
$ws.Range("H4").Value = 4.685748
$ws.Range("I4").Value = 0.02125955722027638
$ws.Range("J4").Value = 0.02125955722027638
$ws.Range("M4").Value = 21.46453166666667
$ws.Range("N4").Value = 64.393595
$ws.Range("O4").Value = 0.3170032929137071
$ws.Range("P4").Value = 0.317003292913707
$ws.Range("Q4").Value = 33.52579544267334
$ws.Range("R4").Value = 301.73215898406
$ws.Range("S4").Value = 0.006739349644714988
$ws.Range("T4").Value = 0.006739349644714988
$ws.Range("G5").Value = 1.561916
$ws.Range("H5").Value = 4.685748
$ws.Range("I5").Value = 0.02125955722027638
$ws.Range("J5").Value = 0.02125955722027638
$ws.Range("M5").Value = 0.2486213333333333
$ws.Range("N5").Value = 0.745864
$ws.Range("O5").Value = 0.003671814627926724
$ws.Range("P5").Value = 0.003671814627926724
$ws.Range("Q5").Value = 0.3883256384746667
$ws.Range("R5").Value = 3.494930746272
$ws.Range("S5").Value = 0.00007806115318465601
$ws.Range("T5").Value = 0.00007806115318465601
$ws.Range("I6").Value = 0.4717415390260894
$ws.Range("J6").Value = 0.4717415390260893
$ws.Range("M6").Value = 2.733663333333333
$ws.Range("N6").Value = 8.200989999999999
$ws.Range("O6").Value = 0.04037266183309663
$ws.Range("P6").Value = 0.04037266183309663
$ws.Range("Q6").Value = 94.74419406969221
$ws.Range("R6").Value = 852.69774662723
$ws.Range("S6").Value = 0.01904546162772486
$ws.Range("T6").Value = 0.01904546162772486
$ws.Range("I7").Value = 0.4717415390260894
$ws.Range("J7").Value = 0.4717415390260893
$ws.Range("O7").Value = 0.6389522306252696
$ws.Range("P7").Value = 0.6389522306252696
$ws.Range("S7").Value = 0.3014203086393175
$ws.Range("T7").Value = 0.3014203086393175
$ws.Range("I8").Value = 0.4717415390260894
$ws.Range("J8").Value = 0.4717415390260893
$ws.Range("M8").Value = 21.46453166666667
$ws.Range("N8").Value = 64.393595
$ws.Range("O8").Value = 0.3170032929137071
$ws.Range("P8").Value = 0.317003292913707
$ws.Range("Q8").Value = 743.9247287858129
$ws.Range("R8").Value = 6695.322559072316
$ws.Range("S8").Value = 0.1495436212754504
$ws.Range("T8").Value = 0.1495436212754504
$ws.Range("I9").Value = 0.4717415390260894
$ws.Range("J9").Value = 0.4717415390260893
$ws.Range("M9").Value = 0.2486213333333333
$ws.Range("N9").Value = 0.745864
$ws.Range("O9").Value = 0.003671814627926724
$ws.Range("P9").Value = 0.003671814627926724
$ws.Range("Q9").Value = 8.616799138347556
$ws.Range("R9").Value = 77.55119224512801
$ws.Range("S9").Value = 0.001732147483596661
$ws.Range("T9").Value = 0.001732147483596661
$ws.Range("G10").Value = 37.00419333333334
$ws.Range("H10").Value = 111.01258
$ws.Range("I10").Value = 0.503671622264046
$ws.Range("J10").Value = 0.5036716222640459
$ws.Range("M10").Value = 2.733663333333333
$ws.Range("N10").Value = 8.200989999999999
$ws.Range("O10").Value = 0.04037266183309663
$ws.Range("P10").Value = 0.04037266183309663
$ws.Range("Q10").Value = 101.1570064949111
$ws.Range("R10").Value = 910.4130584542
$ws.Range("S10").Value = 0.02033456408059351
$ws.Range("T10").Value = 0.02033456408059351
$ws.Range("G11").Value = 37.00419333333334
$ws.Range("H11").Value = 111.01258
$ws.Range("I11").Value = 0.503671622264046
$ws.Range("J11").Value = 0.5036716222640459
$ws.Range("O11").Value = 0.6389522306252696
$ws.Range("P11").Value = 0.6389522306252696
$ws.Range("Q11").Value = 1600.947076774423
$ws.Range("R11").Value = 14408.5236909698
$ws.Range("S11").Value = 0.3218221065482604
$ws.Range("T11").Value = 0.3218221065482604
$ws.Range("G12").Value = 37.00419333333334
$ws.Range("H12").Value = 111.01258
$ws.Range("I12").Value = 0.503671622264046
$ws.Range("J12").Value = 0.5036716222640459
$ws.Range("M12").Value = 21.46453166666667
$ws.Range("N12").Value = 64.393595
$ws.Range("O12").Value = 0.3170032929137071
$ws.Range("P12").Value = 0.317003292913707
$ws.Range("Q12").Value = 794.2776796027891
$ws.Range("R12").Value = 7148.499116425101
$ws.Range("S12").Value = 0.1596655628048914
$ws.Range("T12").Value = 0.1596655628048914
$ws.Range("G13").Value = 37.00419333333334
$ws.Range("H13").Value = 111.01258
$ws.Range("I13").Value = 0.503671622264046
$ws.Range("J13").Value = 0.5036716222640459
$ws.Range("M13").Value = 0.2486213333333333
$ws.Range("N13").Value = 0.745864
$ws.Range("O13").Value = 0.003671814627926724
$ws.Range("P13").Value = 0.003671814627926724
$ws.Range("Q13").Value = 9.200031885457779
$ws.Range("R13").Value = 82.80028696912001
$ws.Range("S13").Value = 0.001849388830300708
$ws.Range("T13").Value = 0.001849388830300708
$ws.Range("G14").Value = 0.2444516666666667
$ws.Range("H14").Value = 0.733355
$ws.Range("I14").Value = 0.003327281489588382
$ws.Range("J14").Value = 0.003327281489588382
$ws.Range("M14").Value = 2.733663333333333
$ws.Range("N14").Value = 8.200989999999999
$ws.Range("O14").Value = 0.04037266183309663
$ws.Range("P14").Value = 0.04037266183309663
$ws.Range("Q14").Value = 0.6682485579388887
$ws.Range("R14").Value = 6.01423702145
$ws.Range("S14").Value = 0.0001343312104026738
$ws.Range("T14").Value = 0.0001343312104026738
$ws.Range("G15").Value = 0.2444516666666667
$ws.Range("H15").Value = 0.733355
$ws.Range("I15").Value = 0.003327281489588382
$ws.Range("J15").Value = 0.003327281489588382
$ws.Range("O15").Value = 0.6389522306252696
$ws.Range("P15").Value = 0.6389522306252696
$ws.Range("Q15").Value = 10.57594142472778
$ws.Range("R15").Value = 95.18347282255
$ws.Range("S15").Value = 0.002125973929690666
$ws.Range("T15").Value = 0.002125973929690666
$ws.Range("G16").Value = 0.2444516666666667
$ws.Range("H16").Value = 0.733355
$ws.Range("I16").Value = 0.003327281489588382
$ws.Range("J16").Value = 0.003327281489588382
$ws.Range("M16").Value = 21.46453166666667
$ws.Range("N16").Value = 64.393595
$ws.Range("O16").Value = 0.3170032929137071
$ws.Range("P16").Value = 0.317003292913707
$ws.Range("Q16").Value = 5.247040540136111
$ws.Range("R16").Value = 47.223364861225
$ws.Range("S16").Value = 0.001054759188650341
$ws.Range("T16").Value = 0.001054759188650341
$ws.Range("G17").Value = 0.2444516666666667
$ws.Range("H17").Value = 0.733355
$ws.Range("I17").Value = 0.003327281489588382
$ws.Range("J17").Value = 0.003327281489588382
$ws.Range("M17").Value = 0.2486213333333333
$ws.Range("N17").Value = 0.745864
$ws.Range("O17").Value = 0.003671814627926724
$ws.Range("P17").Value = 0.003671814627926724
$ws.Range("Q17").Value = 0.06077589930222222
$ws.Range("R17").Value = 0.54698309372
$ws.Range("S17").Value = 0.00001221716084470044
$ws.Range("T17").Value = 0.00001221716084470044
